$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.61
$ws.Range("G2").Value = 1.69
$ws.Range("I2").Value = 7
$ws.Range("L2").Value = 1.35
$ws.Range("P2").Value = 1.89
$ws.Range("Q2").Value = 1.96
$ws.Range("R2").Value = 1.34
$ws.Range("S2").Value = 3.4
$ws.Range("T2").Value = 1.97
$ws.Range("U2").Value = 1.9
$ws.Range("V2").Value = 1.17
$ws.Range("W2").Value = 2.44
$ws.Range("AA2").Value = 210
$ws.Range("AI2").Value = 110
$ws.Range("AM2").Value = 160
$ws.Range("AO2").Value = 160
$ws.Range("F3").Value = 1.24
$ws.Range("N3").Value = 7.8
$ws.Range("P3").Value = 3.7
$ws.Range("R3").Value = 2.1
$ws.Range("T3").Value = 1.76
$ws.Range("U3").Value = 2.12
$ws.Range("X3").Value = 55
$ws.Range("Z3").Value = 160
$ws.Range("AA3").Value = 470
$ws.Range("AE3").Value = 180
$ws.Range("AJ3").Value = 14
$ws.Range("AK3").Value = 16
$ws.Range("AM3").Value = 120
$ws.Range("AN3").Value = 2.98
$ws.Range("AO3").Value = 150
$ws.Range("F4").Value = 3.9
$ws.Range("L4").Value = 1.5
$ws.Range("V4").Value = 1.78
$ws.Range("W4").Value = 1.33
$ws.Range("X4").Value = 9.800000000000001
$ws.Range("Z4").Value = 12.5
$ws.Range("AA4").Value = 29
$ws.Range("AC4").Value = 7.2
$ws.Range("AF4").Value = 26
$ws.Range("AG4").Value = 16
$ws.Range("AJ4").Value = 95
$ws.Range("AK4").Value = 55
$ws.Range("AM4").Value = 170
$ws.Range("AN4").Value = 70
$ws.Range("AO4").Value = 25
$ws.Range("O5").Value = 1.43
$ws.Range("P5").Value = 1.72
$ws.Range("Q5").Value = 2.3
$ws.Range("AE5").Value = 1000
$ws.Range("P6").Value = 1.81
$ws.Range("S6").Value = 3.95
$ws.Range("AH6").Value = 19.5
$ws.Range("AI6").Value = 48
